# Auto-generated edits applying the Sheets workbook diff (scheduled runner update).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 114
$ws.Range("I5").Value = 114
$ws.Range("K5").Value = 114
$ws.Range("M5").Value = 1
$ws.Range("H11").Value = 30484928
$ws.Range("I11").Value = 30484928
$ws.Range("K11").Value = 30484928
$ws.Range("M11").Value = -30484788
$ws.Range("H129").Value = 3975.44
$ws.Range("I129").Value = 2219.6667
$ws.Range("J129").Value = 4963.0625
$ws.Range("K129").Value = 6659.000100000001
$ws.Range("L129").Value = 14889.1875
$ws.Range("M129").Value = -1659.000100000001
$ws.Range("N129").Value = -24889.1875
$ws.Range("H132").Value = 6168.974
$ws.Range("I132").Value = 4141.9
$ws.Range("K132").Value = 12425.7
$ws.Range("M132").Value = -9895.699999999999
$ws.Range("H138").Value = 4785.6616
$ws.Range("I138").Value = 7122.684
$ws.Range("J138").Value = 3879.4695
$ws.Range("K138").Value = 21368.052
$ws.Range("L138").Value = 11638.4085
$ws.Range("M138").Value = -16228.052
$ws.Range("N138").Value = -21918.4085

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 3000
$ws.Range("I26").Value = 3000
$ws.Range("K26").Value = 3000
$ws.Range("M26").Value = -2670
$ws.Range("H130").Value = 20000
$ws.Range("J130").Value = 20000
$ws.Range("L130").Value = 20000
$ws.Range("N130").Value = -30040

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H22").Value = 167449.5
$ws.Range("I22").Value = 167449.5
$ws.Range("K22").Value = 167449.5
$ws.Range("M22").Value = -167276.5
$ws.Range("H29").Value = 15953.5
$ws.Range("I29").Value = 20000
$ws.Range("J29").Value = 14604.667
$ws.Range("K29").Value = 20000
$ws.Range("L29").Value = 14604.667
$ws.Range("M29").Value = -19711
$ws.Range("N29").Value = -15182.667
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H86").Value = 2929.4167
$ws.Range("I86").Value = 2100.3333
$ws.Range("J86").Value = 3758.5
$ws.Range("K86").Value = 2100.3333
$ws.Range("L86").Value = 3758.5
$ws.Range("M86").Value = -977.3332999999998
$ws.Range("N86").Value = -6004.5
$ws.Range("H89").Value = 2929.4167
$ws.Range("I89").Value = 2100.3333
$ws.Range("J89").Value = 3758.5
$ws.Range("K89").Value = 10501.6665
$ws.Range("L89").Value = 18792.5
$ws.Range("M89").Value = -4885.666499999999
$ws.Range("N89").Value = -30024.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 5000
$ws.Range("I26").Value = 5000
$ws.Range("K26").Value = 5000
$ws.Range("M26").Value = -4713
$ws.Range("H31").Value = 3315.7805
$ws.Range("I31").Value = 2164.7666
$ws.Range("K31").Value = 2164.7666
$ws.Range("M31").Value = -1869.7666
$ws.Range("H34").Value = 3315.7805
$ws.Range("I34").Value = 2164.7666
$ws.Range("K34").Value = 2164.7666
$ws.Range("M34").Value = -1962.7666
$ws.Range("H58").Value = 3540.68
$ws.Range("I58").Value = 3462.0476
$ws.Range("K58").Value = 3462.0476
$ws.Range("M58").Value = -3259.0476
$ws.Range("H86").Value = 21955.709
$ws.Range("I86").Value = 22901.904
$ws.Range("J86").Value = 15332.333
$ws.Range("K86").Value = 22901.904
$ws.Range("L86").Value = 15332.333
$ws.Range("M86").Value = -21778.904
$ws.Range("N86").Value = -17578.333
$ws.Range("H89").Value = 21955.709
$ws.Range("I89").Value = 22901.904
$ws.Range("J89").Value = 15332.333
$ws.Range("K89").Value = 114509.52
$ws.Range("L89").Value = 76661.66500000001
$ws.Range("M89").Value = -108893.52
$ws.Range("N89").Value = -87893.66500000001
$ws.Range("H136").Value = 3540.68
$ws.Range("I136").Value = 3462.0476
$ws.Range("K136").Value = 10386.1428
$ws.Range("M136").Value = -7836.1428

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1081.8667
$ws.Range("I11").Value = 1087.7142
$ws.Range("K11").Value = 3263.1426
$ws.Range("M11").Value = -3123.1426
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H25").Value = 4933
$ws.Range("J25").Value = 10999
$ws.Range("L25").Value = 32997
$ws.Range("N25").Value = -33335
$ws.Range("H30").Value = 4933
$ws.Range("J30").Value = 10999
$ws.Range("L30").Value = 32997
$ws.Range("N30").Value = -33201
$ws.Range("H69").Value = 977.5
$ws.Range("J69").Value = 973
$ws.Range("L69").Value = 2919
$ws.Range("N69").Value = -4541
$ws.Range("H72").Value = 977.5
$ws.Range("J72").Value = 973
$ws.Range("L72").Value = 8757
$ws.Range("N72").Value = -16869
$ws.Range("H80").Value = 3866.375
$ws.Range("J80").Value = 4488
$ws.Range("L80").Value = 13464
$ws.Range("N80").Value = -15336
$ws.Range("H83").Value = 3866.375
$ws.Range("J83").Value = 4488
$ws.Range("L83").Value = 40392
$ws.Range("N83").Value = -49752

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H7").Value = 721850.6
$ws.Range("I7").Value = 776862.25
$ws.Range("J7").Value = 6700
$ws.Range("K7").Value = 776862.25
$ws.Range("L7").Value = 6700
$ws.Range("M7").Value = -776750.25
$ws.Range("N7").Value = -6924
$ws.Range("H16").Value = 989.2857
$ws.Range("I16").Value = 900.1111
$ws.Range("J16").Value = 1149.8
$ws.Range("K16").Value = 900.1111
$ws.Range("L16").Value = 1149.8
$ws.Range("M16").Value = -730.1111
$ws.Range("N16").Value = -1489.8
$ws.Range("H40").Value = 2966.2856
$ws.Range("I40").Value = 2175.4375
$ws.Range("J40").Value = 5497
$ws.Range("K40").Value = 2175.4375
$ws.Range("L40").Value = 5497
$ws.Range("M40").Value = -2039.4375
$ws.Range("N40").Value = -5769
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H82").Value = 2146.3157
$ws.Range("J82").Value = 3011.762
$ws.Range("L82").Value = 3011.762
$ws.Range("N82").Value = -3733.762
$ws.Range("H85").Value = 2146.3157
$ws.Range("J85").Value = 3011.762
$ws.Range("L85").Value = 3011.762
$ws.Range("N85").Value = -5507.762000000001
$ws.Range("H126").Value = 721850.6
$ws.Range("I126").Value = 776862.25
$ws.Range("J126").Value = 6700
$ws.Range("K126").Value = 2330586.75
$ws.Range("L126").Value = 20100
$ws.Range("M126").Value = -2328116.75
$ws.Range("N126").Value = -25040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 10874.5
$ws.Range("I32").Value = 10874.5
$ws.Range("K32").Value = 10874.5
$ws.Range("M32").Value = -10557.5
